$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / recalculated mean
$ws.Range("F2").Value = -2
$ws.Range("F6").Value = 3
$ws.Range("F23").Value = -4
$ws.Range("F27").Value = 4
$ws.Range("F33").Value = -1
$ws.Range("F34").Value = -2
$ws.Range("F41").Value = -2
$ws.Range("F44").Value = 2
$ws.Range("F46").Value = 4
$ws.Range("F47").Value = 3
$ws.Range("F52").Value = 0
$ws.Range("F56").Value = -7
